# Update column F (dSF) values on Sheet1 per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 0
    3  = 0
    4  = 2
    5  = -1
    6  = -1
    7  = 2
    8  = 3
    9  = 2
    10 = -3
    11 = 4
    12 = 2
    14 = -2
    15 = -4
    16 = -2
    17 = 5
    18 = 2
    19 = -4
    20 = -4
    21 = -3
    22 = 3
    23 = 2
    24 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
